$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1401.5    # Current Capital
$wsSummary.Range("B4").Value = 1.29      # Total P&L $
$wsSummary.Range("B6").Value = 136       # Total Trades
$wsSummary.Range("B7").Value = 59        # Winning Trades
$wsSummary.Range("B9").Value = 43.38     # Win Rate %

# ---------------------------------------------------------------
# Strategy Status sheet - MarketMaking row (row 5)
# ---------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C5").Value = 101.5
$wsStatus.Range("D5").Value = 103
$wsStatus.Range("E5").Value = 1.18
$wsStatus.Range("F5").Value = 1.5
$wsStatus.Range("G5").Value = 43.69

# ---------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")

# Update existing trade row 137 (Trade #136) - it closes now
$wsAll.Cells.Item(137, 7).Value = 0.088251
$wsAll.Cells.Item(137, 8).Value = "CLOSED"
$wsAll.Cells.Item(137, 9).Value = 10.3136
$wsAll.Cells.Item(137, 10).Value = 0.01
$wsAll.Cells.Item(137, 11).Value = 101.5
$wsAll.Cells.Item(137, 12).Value = "early_exit"
$wsAll.Cells.Item(137, 13).Value = 0.13

# Append new trade row 170 (Trade #169)
$wsAll.Cells.Item(170, 1).Value = 169
$wsAll.Cells.Item(170, 2).NumberFormat = "@"
$wsAll.Cells.Item(170, 2).Value = "2026-02-17"
$wsAll.Cells.Item(170, 3).NumberFormat = "@"
$wsAll.Cells.Item(170, 3).Value = "21:32:44"
$wsAll.Cells.Item(170, 4).Value = "MarketMaking"
$wsAll.Cells.Item(170, 5).Value = "UP"
$wsAll.Cells.Item(170, 6).Value = 0.08
$wsAll.Cells.Item(170, 8).Value = "OPEN"
$wsAll.Cells.Item(170, 9).Value = 0
$wsAll.Cells.Item(170, 10).Value = 0
$wsAll.Cells.Item(170, 11).Value = 101.4941758035408
$wsAll.Cells.Item(170, 13).Value = 0
$wsAll.Cells.Item(170, 14).Value = 0
$wsAll.Cells.Item(170, 15).Value = 0
$wsAll.Cells.Item(170, 16).Value = 0.6
$wsAll.Cells.Item(170, 17).Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------
# MarketMaking sheet (per-strategy mirror of All Trades)
# ---------------------------------------------------------------
$wsMM = $wb.Worksheets.Item("MarketMaking")

# Update existing trade row 104 (Trade #136) - it closes now
$wsMM.Cells.Item(104, 7).Value = 0.088251
$wsMM.Cells.Item(104, 8).Value = "CLOSED"
$wsMM.Cells.Item(104, 9).Value = 10.3136
$wsMM.Cells.Item(104, 10).Value = 0.01
$wsMM.Cells.Item(104, 11).Value = 101.5
$wsMM.Cells.Item(104, 16).Value = "early_exit"
$wsMM.Cells.Item(104, 17).Value = 0.13

# Append new trade row 137 (Trade #169)
$wsMM.Cells.Item(137, 1).Value = 169
$wsMM.Cells.Item(137, 2).NumberFormat = "@"
$wsMM.Cells.Item(137, 2).Value = "2026-02-17"
$wsMM.Cells.Item(137, 3).NumberFormat = "@"
$wsMM.Cells.Item(137, 3).Value = "21:32:44"
$wsMM.Cells.Item(137, 4).Value = "MarketMaking"
$wsMM.Cells.Item(137, 5).Value = "UP"
$wsMM.Cells.Item(137, 6).Value = 0.08
$wsMM.Cells.Item(137, 8).Value = "OPEN"
$wsMM.Cells.Item(137, 9).Value = 0
$wsMM.Cells.Item(137, 10).Value = 0
$wsMM.Cells.Item(137, 11).Value = 101.4941758035408
$wsMM.Cells.Item(137, 12).Value = 0
$wsMM.Cells.Item(137, 13).Value = 0
$wsMM.Cells.Item(137, 14).Value = 0.6
$wsMM.Cells.Item(137, 15).Value = "Normal spread capture: 19600 bps"
$wsMM.Cells.Item(137, 17).Value = 0
